# Insert a new data row (weekly update) at row 20, pushing the
# existing rows 20-56 down to 21-57.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new week's record.
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 44622
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100103
$ws.Cells.Item(20, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(20, 9).Value = 100103002
$ws.Cells.Item(20, 10).Value = "Ciruela"
$ws.Cells.Item(20, 11).Value = "Black Amber"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 270
$ws.Cells.Item(20, 14).Value = 8000
$ws.Cells.Item(20, 15).Value = 8500
$ws.Cells.Item(20, 16).Value = 8222
$ws.Cells.Item(20, 17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 514
$ws.Cells.Item(20, 20).Value = 16
